$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (date serial, new positives, 7-day rolling sum, per 100k inhabitants)
$rows = @(
    @(44330, 0, 0, 0),
    @(44331, 0, 0, 0),
    @(44332, 0, 0, 0),
    @(44333, 1, 1, 48.07692307692308),
    @(44334, 1, 2, 96.15384615384616),
    @(44335, 0, 2, 96.15384615384616),
    @(44336, 0, 2, 96.15384615384616),
    @(44337, 0, 2, 96.15384615384616),
    @(44338, 0, 2, 96.15384615384616),
    @(44339, 0, 2, 96.15384615384616),
    @(44340, 0, 1, 48.07692307692308),
    @(44341, 0, 0, 0),
    @(44342, 0, 0, 0),
    @(44343, 0, 0, 0)
)

$startRow = 256
$endRow = $startRow + $rows.Length - 1

# Copy formatting (style) of the last existing data row down across the new rows first
$ws.Range("A255").Copy()
$ws.Range("A${startRow}:A${endRow}").PasteSpecial(-4122)  # xlPasteFormats

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
